# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    and fill it with the per-fund holdings detail for that quarter (same
#    layout as the other quarterly sheets).
# 2. Insert a new leading row into the "总计" summary sheet for 2022-Q1,
#    pushing the existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value into a cell without Excel's automatic
# number coercion (e.g. "20.94" must stay text, matching the workbook's
# existing convention for these columns), and without leaving a
# permanent number-format override behind on the cell.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# =======================================================================
# Step 1: create the "2022-Q1" sheet, positioned just before "总计"
# =======================================================================
# NOTE: the worksheet reference passed as the "insert before" argument to
# Worksheets.Add() gets rebound to the newly-created sheet once the insert
# happens (it resolves by slot, not by identity) - so re-fetch "总计" by
# name afterwards rather than reusing the pre-Add handle.
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计")

# Borrow the existing header / index-column formatting (bold + border,
# centered) from the prior quarter sheet so the new sheet's style
# exactly matches the rest of the workbook instead of inventing a new one.
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (0-based row counter)
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3
$q1.Range("A6").Value = 4

# Row 2: 010709 安信医药健康主题股票A
Set-TextValue $q1.Range("B2") "010709"
Set-TextValue $q1.Range("C2") "安信医药健康主题股票A"
Set-TextValue $q1.Range("D2") "20.94"
Set-TextValue $q1.Range("E2") "89.03"
Set-TextValue $q1.Range("F2") "4.52"
Set-TextValue $q1.Range("G2") "0.9465"
$q1.Range("H2").Value = 4

# Row 3: 010710 安信医药健康主题股票C
Set-TextValue $q1.Range("B3") "010710"
Set-TextValue $q1.Range("C3") "安信医药健康主题股票C"
Set-TextValue $q1.Range("D3") "10.09"
Set-TextValue $q1.Range("E3") "89.03"
Set-TextValue $q1.Range("F3") "4.52"
Set-TextValue $q1.Range("G3") "0.4561"
$q1.Range("H3").Value = 4

# Row 4: 003516 国泰融安多策略灵活配置混合
Set-TextValue $q1.Range("B4") "003516"
Set-TextValue $q1.Range("C4") "国泰融安多策略灵活配置混合"
Set-TextValue $q1.Range("D4") "11.18"
Set-TextValue $q1.Range("E4") "71.30"
Set-TextValue $q1.Range("F4") "1.89"
Set-TextValue $q1.Range("G4") "0.2113"
$q1.Range("H4").Value = 7

# Row 5: 020023 国泰事件驱动策略混合
Set-TextValue $q1.Range("B5") "020023"
Set-TextValue $q1.Range("C5") "国泰事件驱动策略混合"
Set-TextValue $q1.Range("D5") "3.34"
Set-TextValue $q1.Range("E5") "74.02"
Set-TextValue $q1.Range("F5") "1.74"
Set-TextValue $q1.Range("G5") "0.0581"
$q1.Range("H5").Value = 8

# Row 6: 217021 招商优势企业混合
Set-TextValue $q1.Range("B6") "217021"
Set-TextValue $q1.Range("C6") "招商优势企业混合"
Set-TextValue $q1.Range("D6") "0.36"
Set-TextValue $q1.Range("E6") "69.72"
Set-TextValue $q1.Range("F6") "4.76"
Set-TextValue $q1.Range("G6") "0.0171"
$q1.Range("H6").Value = 8

# =======================================================================
# Step 2: prepend a "2022-Q1" row to the "总计" summary sheet, shifting
# the previously-existing rows down by one and renumbering column A.
# =======================================================================

# Grow the table by one row, copying formatting from the existing last
# row so the new row 7 matches the others' styling exactly.
$total.Range("A6:D6").Copy()
$total.Range("A7:D7").PasteSpecial(-4122)

# Shift existing data (old rows 2-6) down to rows 3-7.
# NOTE: `.Value` as a plain getter misbehaves in this host (it stringifies
# the property descriptor instead of returning the cell's contents) - use
# `.Value2` to read, which works correctly; `.Value` remains fine as a
# setter (used elsewhere in this script).
$total.Range("B7").Value = $total.Range("B6").Value2
$total.Range("C7").Value = $total.Range("C6").Value2
$total.Range("D7").Value = $total.Range("D6").Value2

$total.Range("B6").Value = $total.Range("B5").Value2
$total.Range("C6").Value = $total.Range("C5").Value2
$total.Range("D6").Value = $total.Range("D5").Value2

$total.Range("B5").Value = $total.Range("B4").Value2
$total.Range("C5").Value = $total.Range("C4").Value2
$total.Range("D5").Value = $total.Range("D4").Value2

$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

# New first data row: 2022-Q1
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.69

# Renumber the 0-based index column for all 6 data rows.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
